$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D ("batsman") to make room for
# "ownTeam" and "oppTeam"
$ws.Range("D1:E1").EntireColumn.Insert()

# Set new header values
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Set new row 2 values
$ws.Range("D2").Value = "Royal Challengers Bangalore"
$ws.Range("E2").Value = "Delhi Capitals"
